# Prevent division by zero error in city_economic_calc.py for CHP usage.
#
# The three buildings' net floor areas (tab_ease_building_net_floor_area,
# column E) are bumped up, and the stale, manually pre-computed heat/
# electricity-demand figures that depended on the old areas (columns H/I,
# which downstream code was dividing by / using inconsistently and could
# hit a ZeroDivisionError for CHP sizing) are cleared out so they get
# recomputed rather than read as stale inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 2 (building id = 1) ---
$ws.Range("E2").Value = 200          # net floor area: 100 -> 200
$ws.Range("H2").ClearContents()      # drop stale "=E2*150" heat-demand formula
$ws.Range("I2").ClearContents()      # drop stale electr. demand value (2500)
$ws.Range("J2").Value = 0            # usable pv roof area now explicit 0

# --- Row 3 (building id = 2) ---
$ws.Range("E3").Value = 400          # net floor area: 150 -> 400
$ws.Range("L3").Value = 8            # total number of occupants: 5 -> 8

# --- Row 4 (building id = 3) ---
$ws.Range("E4").Value = 300          # net floor area: 125 -> 300
$ws.Range("H4").ClearContents()      # drop stale heat-demand value (10000)
$ws.Range("I4").ClearContents()      # drop stale electr. demand value (3200)

# --- View state: leave the saved selection on the last-edited cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E4").Select()
